$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-12-10 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-12-11 Wednesday", 2) | Out-Null

# Update the division-problem table. The table has 5 "answer" rows (1, 5, 9,
# 13, 17) interleaved with blank rows, each with 5 columns. Some values are
# duplicated within the table (e.g. "15÷2=7, 1" appears twice), so cells are
# addressed directly by row/column rather than via Find/Replace to avoid
# accidentally changing the wrong occurrence.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "79÷7=11, 2"
$t.Cell(1,2).Range.Text  = "56÷8=7, 0"
$t.Cell(1,3).Range.Text  = "53÷6=8, 5"
$t.Cell(1,4).Range.Text  = "11÷3=3, 2"
$t.Cell(1,5).Range.Text  = "75÷5=15, 0"

$t.Cell(5,1).Range.Text  = "69÷4=17, 1"
$t.Cell(5,2).Range.Text  = "56÷6=9, 2"
$t.Cell(5,3).Range.Text  = "94÷8=11, 6"
$t.Cell(5,4).Range.Text  = "21÷5=4, 1"
$t.Cell(5,5).Range.Text  = "84÷2=42, 0"

$t.Cell(9,1).Range.Text  = "74÷9=8, 2"
$t.Cell(9,2).Range.Text  = "12÷5=2, 2"
$t.Cell(9,3).Range.Text  = "69÷3=23, 0"
$t.Cell(9,4).Range.Text  = "97÷6=16, 1"
$t.Cell(9,5).Range.Text  = "38÷6=6, 2"

$t.Cell(13,1).Range.Text = "42÷9=4, 6"
$t.Cell(13,2).Range.Text = "57÷9=6, 3"
$t.Cell(13,3).Range.Text = "94÷4=23, 2"
$t.Cell(13,4).Range.Text = "85÷8=10, 5"
$t.Cell(13,5).Range.Text = "37÷4=9, 1"

$t.Cell(17,1).Range.Text = "16÷7=2, 2"
$t.Cell(17,2).Range.Text = "27÷3=9, 0"
$t.Cell(17,3).Range.Text = "56÷7=8, 0"
$t.Cell(17,4).Range.Text = "36÷2=18, 0"
$t.Cell(17,5).Range.Text = "73÷5=14, 3"
